$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same table layout with
# updated "想去人数" (want-to-go count) values in column F.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5412
    $ws.Range("F3").Value = 165
    $ws.Range("F4").Value = 930
}
